$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row holding "4.4 Retorna ao passo 5" (old row 18) - its content is
# superseded by the edit and the row is dropped entirely, shifting subsequent
# rows up by one.
$ws.Rows.Item(18).Delete()

# Rename the first "System response" entry: the generic step is replaced with
# the more specific wording for accessory components.
$ws.Range("D7").Value = "1. Apresenta lista de componentes acessórios"
